$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J (match style of existing header H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for columns I and J, rows 2-15
$data = @(
    @(8, 8),
    @(7, 7),
    @(5, 5),
    @(9, 9),
    @(9, 9),
    @(6, 7),
    @(5, 5),
    @(8, 8),
    @(7, 8),
    @(5, 6),
    @(5, 5),
    @(5, 6),
    @(6, 6),
    @(8, 8)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
